# Hortaliza, Femacal de La Calera - Poroto verde
# A new weekly price record is inserted as row 471 (Limache / $/saco 25 kilos),
# pushing every existing record from the old row 471 down through row 537
# into rows 472-538. The sheet's used range grows from A1:R537 to A1:R538.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 471; Excel shifts rows 471:537 down to 472:538
# and extends the dimension/used range automatically.
$ws.Rows.Item(471).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A471").Value = 3
$ws.Range("B471").Value = "Femacal de La Calera"
$ws.Range("C471").Value = "Coquimbo"
$ws.Range("D471").Value = 44984
$ws.Range("E471").Value = 5
$ws.Range("F471").Value = 100112031
$ws.Range("G471").Value = "Poroto verde"
$ws.Range("H471").Value = "Magnum"
$ws.Range("I471").Value = "Primera"
$ws.Range("J471").Value = 78
$ws.Range("K471").Value = 25000
$ws.Range("L471").Value = 26000
$ws.Range("M471").Value = 25487
$ws.Range("N471").Value = "$/saco 25 kilos"
$ws.Range("O471").Value = "Limache"
$ws.Range("P471").Value = 1019
$ws.Range("Q471").Value = 25
$ws.Range("R471").Value = "Hortaliza"
